# Updated cryptos list on Mon May 20 20:34:58 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.573.59"
$ws.Range("E2").Value = "  +4.96%  "
$ws.Range("D3").Value = "3.445.61"
$ws.Range("E3").Value = "  +11.89%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'586.48"
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("D6").Value = "'185.26"
$ws.Range("E6").Value = "  +8.95%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.447.61"
$ws.Range("E8").Value = "  +12.04%  "
$ws.Range("D9").Value = "'0.532"
$ws.Range("E9").Value = "  +4.43%  "
$ws.Range("D10").Value = "'6.58"
$ws.Range("E10").Value = "  +3.97%  "
$ws.Range("D11").Value = "'0.157"
$ws.Range("E11").Value = "  +5.38%  "
$ws.Range("D12").Value = "'0.486"
$ws.Range("E12").Value = "  +3.34%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "'38.34"
$ws.Range("E13").Value = "  +7.21%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "'0.0000249"
$ws.Range("E14").Value = "  +4.04%  "
$ws.Range("D15").Value = "4.016.20"
$ws.Range("E15").Value = "  +11.81%  "
$ws.Range("D16").Value = "69.713.10"
$ws.Range("E16").Value = "  +5.29%  "
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").Value = "3.443.28"
$ws.Range("E18").Value = "  +11.79%  "
$ws.Range("D19").Value = "'7.38"
$ws.Range("E19").Value = "  +6.14%  "
$ws.Range("D20").Value = "'16.86"
$ws.Range("E20").Value = "  +1.29%  "
$ws.Range("D21").Value = "'500.53"
$ws.Range("E21").Value = "  +2.65%  "
$ws.Range("D22").Value = "'8.75"
$ws.Range("E22").Value = "  +13.81%  "
$ws.Range("D23").Value = "'0.725"
$ws.Range("E23").Value = "  +5.62%  "
$ws.Range("D24").Value = "'86.36"
$ws.Range("E24").Value = "  +4.63%  "
$ws.Range("D25").Value = "'13.27"
$ws.Range("E25").Value = "  +4.94%  "
$ws.Range("D26").Value = "'2.37"
$ws.Range("E26").Value = "  +7.64%  "
$ws.Range("D27").Value = "'10.73"
$ws.Range("E27").Value = "  +4.83%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  +11.03%  "
$ws.Range("D30").Value = "'8.11"
$ws.Range("E30").Value = "  +2.97%  "
$ws.Range("E31").Value = "  +4.04%  "
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "'0.0000106"
$ws.Range("E32").Value = "  +17.72%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'30.31"
$ws.Range("E33").Value = "  +9.32%  "
$ws.Range("D34").Value = "'0.116"
$ws.Range("E34").Value = "  +4.28%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'6.10"
$ws.Range("E36").Value = "  +9.49%  "
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  +5.28%  "
$ws.Range("D38").Value = "'47.81"
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("E39").Value = "  +9.63%  "
$ws.Range("D40").Value = "'2.10"
$ws.Range("E40").Value = "  +6.94%  "
$ws.Range("E41").Value = "  +4.94%  "
$ws.Range("D42").Value = "'50.15"
$ws.Range("E42").Value = "  +2.20%  "
$ws.Range("D43").Value = "'8.65"
$ws.Range("E43").Value = "  +4.72%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.85"
$ws.Range("E44").Value = "  +13.60%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'409.35"
$ws.Range("E45").Value = "  +12.11%  "
$ws.Range("D46").Value = "2.942.68"
$ws.Range("E46").Value = "  +5.64%  "
$ws.Range("D47").Value = "'28.10"
$ws.Range("E47").Value = "  +15.02%  "
$ws.Range("E48").Value = "  +4.61%  "
$ws.Range("D49").Value = "'134.94"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").Value = "'2.44"
$ws.Range("E51").Value = "  +13.37%  "
